# "weird hangup with put"
#
# Changes:
#  - Data models sheet: Boats.owner type changes from "int " to "string"
#  - Data models sheet: Users entity loses its creation_date and last_modified
#    rows (the shared string "last_modified" becomes unused and drops out of
#    the workbook entirely once no cell references it)
#  - methods sheet: the "/users" column (B) no longer lists PUT or PATCH as
#    supported methods
#  - selection / active-sheet bookkeeping is updated to match: "methods"
#    sheet cursor moves to B10, then "Data models" becomes the active /
#    selected tab with the cursor on A8 (previously "architecture" was active)

$wb = $excel.ActiveWorkbook

$wsModels = $wb.Worksheets.Item("Data models")
$wsMethods = $wb.Worksheets.Item("methods")

# Boats.owner: int -> string
$wsModels.Range("B7").Value = "string"

# Users: drop the creation_date (row 13) and last_modified (row 14) rows
$wsModels.Range("A13:C14").ClearContents()

# methods sheet: /users no longer supports PUT (row 9) or PATCH (row 10)
$wsMethods.Range("B9").ClearContents()
$wsMethods.Range("B10").ClearContents()
$wsMethods.Range("B10").Select()

# Make "Data models" the active / selected sheet with A8 selected
$wsModels.Activate()
$wsModels.Range("A8").Select()
